$wb = $excel.ActiveWorkbook

# Sheet "展览" (sheet1): update F3, F9, F10, F11
$wsExhibit = $wb.Worksheets.Item("展览")
$wsExhibit.Range("F3").Value = 507
$wsExhibit.Range("F9").Value = 129
$wsExhibit.Range("F10").Value = 2293
$wsExhibit.Range("F11").Value = 11

# Sheet "全部类型" (sheet4): update F4, F10, F11, F12
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F4").Value = 507
$wsAll.Range("F10").Value = 129
$wsAll.Range("F11").Value = 2293
$wsAll.Range("F12").Value = 11
